# Weekly price-list update for "Hortaliza, Vega Monumental Concepción - Brócoli".
# A new price observation is inserted as row 432 (Fecha serial 45013 = 2023-03-28),
# pushing the existing rows 432:457 down to 433:458 (dimension grows to A1:R458).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh blank row above the current row 432; this shifts every row
# from 432 through 457 down by one (to 433 through 458) and copies the
# column formatting (e.g. the date number format on column D) down with it.
$ws.Rows(432).Insert()

# Populate the newly inserted row 432 with the new weekly data point.
$ws.Cells.Item(432, 1).Value = 11
$ws.Cells.Item(432, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(432, 3).Value = "Bíobío"
$ws.Cells.Item(432, 4).Value = 45013
$ws.Cells.Item(432, 5).Value = 8
$ws.Cells.Item(432, 6).Value = 100112023
$ws.Cells.Item(432, 7).Value = "Brócoli"
$ws.Cells.Item(432, 8).Value = "Sin especificar"
$ws.Cells.Item(432, 9).Value = "Primera"
$ws.Cells.Item(432, 10).Value = 2200
$ws.Cells.Item(432, 11).Value = 900
$ws.Cells.Item(432, 12).Value = 1000
$ws.Cells.Item(432, 13).Value = 955
$ws.Cells.Item(432, 14).Value = "$/unidad"
$ws.Cells.Item(432, 15).Value = "Región Metropolitana"
$ws.Cells.Item(432, 16).Value = 955
$ws.Cells.Item(432, 17).Value = 1
$ws.Cells.Item(432, 18).Value = "Hortaliza"
